$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.19"
$ws.Range("E2").Value = "'5.38%"
$ws.Range("D3").Value = "'34.82"
$ws.Range("E3").Value = "'12.34%"
$ws.Range("D4").Value = "'5.174"
$ws.Range("E4").Value = "'4.41%"
$ws.Range("D5").Value = "'0.07795"
$ws.Range("E5").Value = "'6.05%"
$ws.Range("D6").Value = "'2.295"
$ws.Range("E6").Value = "'-0.40%"
$ws.Range("D7").Value = "'8.058"
$ws.Range("E7").Value = "'4.35%"
$ws.Range("D8").Value = "'3.992"
$ws.Range("E8").Value = "'6.64%"
$ws.Range("D9").Value = "'0.9236"
$ws.Range("E9").Value = "'1.15%"
$ws.Range("D10").Value = "'0.1006"
$ws.Range("E10").Value = "'7.61%"
$ws.Range("D11").Value = "'0.1832"
$ws.Range("E11").Value = "'7.58%"
$ws.Range("D12").Value = "'0.08551"
$ws.Range("E12").Value = "'3.58%"
$ws.Range("D13").Value = "'0.03394"
$ws.Range("E13").Value = "'9.24%"
$ws.Range("D14").Value = "'0.09907"
$ws.Range("E14").Value = "'-0.80%"
$ws.Range("D15").Value = "'0.001481"
$ws.Range("E15").Value = "'-0.99%"
$ws.Range("D16").Value = "'0.04658"
$ws.Range("E16").Value = "'2.98%"
$ws.Range("D17").Value = "'0.005790"
$ws.Range("E17").Value = "'0.83%"
$ws.Range("E18").Value = "'-0.02%"
$ws.Range("D19").Value = "'2.127"
$ws.Range("E19").Value = "'5.90%"
$ws.Range("E20").Value = "'2.84%"
$ws.Range("D21").Value = "'0.1328"
$ws.Range("E21").Value = "'3.14%"
$ws.Range("D22").Value = "'4.552"
$ws.Range("E22").Value = "'9.63%"
$ws.Range("D23").Value = "'0.2272"
$ws.Range("E23").Value = "'8.24%"
$ws.Range("D24").Value = "'0.001220"
$ws.Range("E24").Value = "'0.74%"
$ws.Range("D25").Value = "'0.004338"
$ws.Range("E25").Value = "'3.91%"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'0.01%"
$ws.Range("D27").Value = "'0.0003401"
$ws.Range("E27").Value = "'0.13%"
$ws.Range("D39").Value = "'0.01743"
$ws.Range("E39").Value = "'11.05%"
$ws.Range("D40").Value = "'0.04737"
$ws.Range("E40").Value = "'5.92%"
$ws.Range("D41").Value = "'0.007717"
$ws.Range("E41").Value = "'4.51%"
$ws.Range("D42").Value = "'0.1412"
$ws.Range("E42").Value = "'6.09%"
$ws.Range("E43").Value = "'-22.43%"
$ws.Range("D44").Value = "'0.002302"
$ws.Range("E44").Value = "'2.24%"
$ws.Range("D45").Value = "'0.009978"
$ws.Range("E45").Value = "'13.90%"
$ws.Range("D46").Value = "'0.00006083"
$ws.Range("E46").Value = "'-0.56%"
$ws.Range("E47").Value = "'0.02%"
$ws.Range("E48").Value = "'51.25%"
$ws.Range("D49").Value = "'0.002692"
$ws.Range("E49").Value = "'34.52%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.02%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.02%"
